# Loan RBI, Variable Instalments
# - Insert a new (blank) column N on the "Repayment Schedule" sheet, pushing
#   the existing "Late"/"Outstanding" columns one to the right.
# - Give the newly inserted column a width matching its neighbour (column M).
# - Make "Repayment Schedule" the active sheet / active tab, with the
#   selection parked on S7 (previously the "Transactions" sheet was active).

$wb = $excel.ActiveWorkbook

$repay = $wb.Worksheets.Item("Repayment Schedule")

# Remember the width of column M before we shift things around, so the new
# column N can be sized the same way.
$colMWidth = $repay.Columns("M").ColumnWidth

# Insert a new blank column at N - everything from the old N ("Late") onward
# shifts one column to the right (N->O, O->P, P->Q, ...).
$repay.Columns("N").Insert() | Out-Null

# Size the freshly inserted column like its left neighbour.
$repay.Columns("N").ColumnWidth = $colMWidth

# Switch the active sheet/tab to "Repayment Schedule" and park the selection
# on S7 (this also clears the "Transactions" sheet's tab-selected state).
$repay.Activate() | Out-Null
$repay.Range("S7").Select() | Out-Null
